# Updated cryptos list on Tue Mar 28 22:58:33 UTC 2023 with GitHub Actions
#
# Refresh the Price (D) and Volume(1h) (E) columns on Sheet1 with the
# latest scraped snapshot values. Several Price cells are numeric-looking
# strings (e.g. "312.79", "0.5222") that must stay stored as TEXT, just
# like the rest of the sheet -- so for those we force NumberFormat to
# "@" (Text) immediately before writing the value, which mirrors how
# Excel keeps a quote-prefixed / text-formatted entry from being
# reinterpreted as a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "27.233.96"
$ws.Cells.Item(2, 5).Value = "  +0.34%  "
$ws.Cells.Item(3, 4).Value = "1.771.94"
$ws.Cells.Item(3, 5).Value = "  +3.51%  "
$ws.Cells.Item(4, 5).Value = "  +0.02%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "312.79"
$ws.Cells.Item(5, 5).Value = "  +1.16%  "
$ws.Cells.Item(6, 5).Value = "  +0.03%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.5222"
$ws.Cells.Item(7, 5).Value = "  +9.54%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.3663"
$ws.Cells.Item(8, 5).Value = "  +6.51%  "
$ws.Cells.Item(9, 5).Value = "  +1.47%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.07354"
$ws.Cells.Item(10, 5).Value = "  +0.98%  "
$ws.Cells.Item(11, 5).Value = "  +4.24%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "1.001"
$ws.Cells.Item(12, 5).Value = "  +0.00%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "20.45"
$ws.Cells.Item(13, 5).Value = "  +3.08%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "6.058"
$ws.Cells.Item(14, 5).Value = "  +3.49%  "
$ws.Cells.Item(15, 4).Value = "1.770.75"
$ws.Cells.Item(15, 5).Value = "  +3.36%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "6.929"
$ws.Cells.Item(16, 5).Value = "  +1.36%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "88.63"
$ws.Cells.Item(17, 5).Value = "  -0.08%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.00001044"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.06438"
$ws.Cells.Item(19, 5).Value = "  +1.25%  "
$ws.Cells.Item(20, 5).Value = "  +0.03%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "16.70"
$ws.Cells.Item(21, 5).Value = "  +1.40%  "
$ws.Cells.Item(22, 5).Value = "  +3.37%  "
$ws.Cells.Item(23, 4).Value = "27.272.93"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "11.23"
$ws.Cells.Item(24, 5).Value = "  +3.90%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.128"
$ws.Cells.Item(25, 5).Value = "  +1.74%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "154.94"
$ws.Cells.Item(26, 5).Value = "  +1.66%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "20.10"
$ws.Cells.Item(27, 5).Value = "  +2.35%  "
$ws.Cells.Item(28, 4).Value = "1.972.79"
$ws.Cells.Item(28, 5).Value = "  +3.40%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "2.311"
$ws.Cells.Item(29, 5).Value = "  +10.80%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "121.03"
$ws.Cells.Item(30, 5).Value = "  +0.74%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "1.053"
$ws.Cells.Item(31, 5).Value = "  +3.90%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.09764"
$ws.Cells.Item(32, 5).Value = "  +5.46%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "5.545"
$ws.Cells.Item(33, 5).Value = "  +4.63%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "3.622"
$ws.Cells.Item(34, 5).Value = "  +0.95%  "
$ws.Cells.Item(35, 5).Value = "  +1.48%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.05944"
$ws.Cells.Item(36, 5).Value = "  +0.83%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "11.16"
$ws.Cells.Item(37, 5).Value = "  +1.05%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "4.819"
$ws.Cells.Item(38, 5).Value = "  +1.36%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.6116"
$ws.Cells.Item(39, 5).Value = "  +3.22%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.2013"
$ws.Cells.Item(40, 5).Value = "  -0.02%  "
$ws.Cells.Item(41, 5).Value = "  +1.27%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "8.055"
$ws.Cells.Item(42, 5).Value = "  +7.33%  "
$ws.Cells.Item(43, 5).Value = "  +1.99%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "13.03"
$ws.Cells.Item(44, 5).Value = "  +2.97%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.5753"
$ws.Cells.Item(45, 5).Value = "  +2.37%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "3.621"
$ws.Cells.Item(46, 5).Value = "  +1.46%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "121.04"
$ws.Cells.Item(47, 5).Value = "  +2.20%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "1.872"
$ws.Cells.Item(48, 5).Value = "  +1.78%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "1.113"
$ws.Cells.Item(49, 5).Value = "  +2.39%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.06699"
$ws.Cells.Item(50, 5).Value = "  +0.99%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.001"
$ws.Cells.Item(51, 5).Value = "  +0.06%  "
